# Update the "Periodo Mora" column (E) values in the worker arrears table.
# The account-statement period codes are renumbered from the previous
# descending order (2507, 2506, 2505, 2504) to an ascending order
# (2504, 2505, 2506, 2507) as part of refreshing the EC database with the
# first batch of new account statements.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2504"
$ws.Range("E17").Value = "2505"
$ws.Range("E18").Value = "2506"
$ws.Range("E19").Value = "2507"
